$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append rows 19-35: a repeat (re-scraped, reshuffled) copy of the
# existing match-log rows for Shreyas Iyer (c), doubling the sheet to A1:K35.
# Every cell in this sheet -- including numeric-looking values like run
# counts and strike rates -- is stored as literal TEXT (not Number), so we
# force a Text number format before writing each value to preserve that.

# Row 19
$ws.Range("A19:K19").NumberFormat = "@"
$ws.Range("A19").Value = " Dubai (DSC)"
$ws.Range("B19").Value = " October 14 2020"
$ws.Range("C19").Value = "Capitals won by 13 runs"
$ws.Range("D19").Value = "Delhi Capitals"
$ws.Range("E19").Value = "Rajasthan Royals"
$ws.Range("F19").Value = "Shreyas Iyer (c)"
$ws.Range("G19").Value = "53"
$ws.Range("H19").Value = "43"
$ws.Range("I19").Value = "3"
$ws.Range("J19").Value = "2"
$ws.Range("K19").Value = "123.25"

# Row 20
$ws.Range("A20:K20").NumberFormat = "@"
$ws.Range("A20").Value = " Abu Dhabi"
$ws.Range("B20").Value = " October 11 2020"
$ws.Range("C20").Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Range("D20").Value = "Delhi Capitals"
$ws.Range("E20").Value = "Mumbai Indians"
$ws.Range("F20").Value = "Shreyas Iyer (c)"
$ws.Range("G20").Value = "42"
$ws.Range("H20").Value = "33"
$ws.Range("I20").Value = "5"
$ws.Range("J20").Value = "0"
$ws.Range("K20").Value = "127.27"

# Row 21
$ws.Range("A21:K21").NumberFormat = "@"
$ws.Range("A21").Value = " Abu Dhabi"
$ws.Range("B21").Value = " November 02 2020"
$ws.Range("C21").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D21").Value = "Delhi Capitals"
$ws.Range("E21").Value = "Royal Challengers Bangalore"
$ws.Range("F21").Value = "Shreyas Iyer (c)"
$ws.Range("G21").Value = "7"
$ws.Range("H21").Value = "9"
$ws.Range("I21").Value = "0"
$ws.Range("J21").Value = "0"
$ws.Range("K21").Value = "77.77"

# Row 22
$ws.Range("A22:K22").NumberFormat = "@"
$ws.Range("A22").Value = " Dubai (DSC)"
$ws.Range("B22").Value = " November 05 2020"
$ws.Range("C22").Value = "Mumbai won by 57 runs"
$ws.Range("D22").Value = "Delhi Capitals"
$ws.Range("E22").Value = "Mumbai Indians"
$ws.Range("F22").Value = "Shreyas Iyer (c)"
$ws.Range("G22").Value = "12"
$ws.Range("H22").Value = "8"
$ws.Range("I22").Value = "3"
$ws.Range("J22").Value = "0"
$ws.Range("K22").Value = "150.00"

# Row 23
$ws.Range("A23:K23").NumberFormat = "@"
$ws.Range("A23").Value = " Dubai (DSC)"
$ws.Range("B23").Value = " October 27 2020"
$ws.Range("C23").Value = "Sunrisers won by 88 runs"
$ws.Range("D23").Value = "Delhi Capitals"
$ws.Range("E23").Value = "Sunrisers Hyderabad"
$ws.Range("F23").Value = "Shreyas Iyer (c)"
$ws.Range("G23").Value = "7"
$ws.Range("H23").Value = "12"
$ws.Range("I23").Value = "0"
$ws.Range("J23").Value = "0"
$ws.Range("K23").Value = "58.33"

# Row 24
$ws.Range("A24:K24").NumberFormat = "@"
$ws.Range("A24").Value = " Abu Dhabi"
$ws.Range("B24").Value = " September 29 2020"
$ws.Range("C24").Value = "Sunrisers won by 15 runs"
$ws.Range("D24").Value = "Delhi Capitals"
$ws.Range("E24").Value = "Sunrisers Hyderabad"
$ws.Range("F24").Value = "Shreyas Iyer (c)"
$ws.Range("G24").Value = "17"
$ws.Range("H24").Value = "21"
$ws.Range("I24").Value = "2"
$ws.Range("J24").Value = "0"
$ws.Range("K24").Value = "80.95"

# Row 25
$ws.Range("A25:K25").NumberFormat = "@"
$ws.Range("A25").Value = " Abu Dhabi"
$ws.Range("B25").Value = " October 24 2020"
$ws.Range("C25").Value = "KKR won by 59 runs"
$ws.Range("D25").Value = "Delhi Capitals"
$ws.Range("E25").Value = "Kolkata Knight Riders"
$ws.Range("F25").Value = "Shreyas Iyer (c)"
$ws.Range("G25").Value = "47"
$ws.Range("H25").Value = "38"
$ws.Range("I25").Value = "5"
$ws.Range("J25").Value = "0"
$ws.Range("K25").Value = "123.68"

# Row 26
$ws.Range("A26:K26").NumberFormat = "@"
$ws.Range("A26").Value = " Sharjah"
$ws.Range("B26").Value = " October 17 2020"
$ws.Range("C26").Value = "Capitals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D26").Value = "Delhi Capitals"
$ws.Range("E26").Value = "Chennai Super Kings"
$ws.Range("F26").Value = "Shreyas Iyer (c)"
$ws.Range("G26").Value = "23"
$ws.Range("H26").Value = "23"
$ws.Range("I26").Value = "1"
$ws.Range("J26").Value = "1"
$ws.Range("K26").Value = "100.00"

# Row 27
$ws.Range("A27:K27").NumberFormat = "@"
$ws.Range("A27").Value = " Dubai (DSC)"
$ws.Range("B27").Value = " October 31 2020"
$ws.Range("C27").Value = "Mumbai won by 9 wickets (with 34 balls remaining)"
$ws.Range("D27").Value = "Delhi Capitals"
$ws.Range("E27").Value = "Mumbai Indians"
$ws.Range("F27").Value = "Shreyas Iyer (c)"
$ws.Range("G27").Value = "25"
$ws.Range("H27").Value = "29"
$ws.Range("I27").Value = "1"
$ws.Range("J27").Value = "1"
$ws.Range("K27").Value = "86.20"

# Row 28
$ws.Range("A28:K28").NumberFormat = "@"
$ws.Range("A28").Value = " Dubai (DSC)"
$ws.Range("B28").Value = " November 10 2020"
$ws.Range("C28").Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Range("D28").Value = "Delhi Capitals"
$ws.Range("E28").Value = "Mumbai Indians"
$ws.Range("F28").Value = "Shreyas Iyer (c)"
$ws.Range("G28").Value = "65"
$ws.Range("H28").Value = "50"
$ws.Range("I28").Value = "6"
$ws.Range("J28").Value = "2"
$ws.Range("K28").Value = "130.00"

# Row 29
$ws.Range("A29:K29").NumberFormat = "@"
$ws.Range("A29").Value = " Dubai (DSC)"
$ws.Range("B29").Value = " October 20 2020"
$ws.Range("C29").Value = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Range("D29").Value = "Delhi Capitals"
$ws.Range("E29").Value = "Kings XI Punjab"
$ws.Range("F29").Value = "Shreyas Iyer (c)"
$ws.Range("G29").Value = "14"
$ws.Range("H29").Value = "12"
$ws.Range("I29").Value = "0"
$ws.Range("J29").Value = "1"
$ws.Range("K29").Value = "116.66"

# Row 30
$ws.Range("A30:K30").NumberFormat = "@"
$ws.Range("A30").Value = " Dubai (DSC)"
$ws.Range("B30").Value = " October 05 2020"
$ws.Range("C30").Value = "Capitals won by 59 runs"
$ws.Range("D30").Value = "Delhi Capitals"
$ws.Range("E30").Value = "Royal Challengers Bangalore"
$ws.Range("F30").Value = "Shreyas Iyer (c)"
$ws.Range("G30").Value = "11"
$ws.Range("H30").Value = "13"
$ws.Range("I30").Value = "1"
$ws.Range("J30").Value = "0"
$ws.Range("K30").Value = "84.61"

# Row 31
$ws.Range("A31:K31").NumberFormat = "@"
$ws.Range("A31").Value = " Dubai (DSC)"
$ws.Range("B31").Value = " September 20 2020"
$ws.Range("C31").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D31").Value = "Delhi Capitals"
$ws.Range("E31").Value = "Kings XI Punjab"
$ws.Range("F31").Value = "Shreyas Iyer (c)"
$ws.Range("G31").Value = "39"
$ws.Range("H31").Value = "32"
$ws.Range("I31").Value = "0"
$ws.Range("J31").Value = "3"
$ws.Range("K31").Value = "121.87"

# Row 32
$ws.Range("A32:K32").NumberFormat = "@"
$ws.Range("A32").Value = " Abu Dhabi"
$ws.Range("B32").Value = " November 08 2020"
$ws.Range("C32").Value = "Capitals won by 17 runs"
$ws.Range("D32").Value = "Delhi Capitals"
$ws.Range("E32").Value = "Sunrisers Hyderabad"
$ws.Range("F32").Value = "Shreyas Iyer (c)"
$ws.Range("G32").Value = "21"
$ws.Range("H32").Value = "20"
$ws.Range("I32").Value = "1"
$ws.Range("J32").Value = "0"
$ws.Range("K32").Value = "105.00"

# Row 33
$ws.Range("A33:K33").NumberFormat = "@"
$ws.Range("A33").Value = " Sharjah"
$ws.Range("B33").Value = " October 03 2020"
$ws.Range("C33").Value = "Capitals won by 18 runs"
$ws.Range("D33").Value = "Delhi Capitals"
$ws.Range("E33").Value = "Kolkata Knight Riders"
$ws.Range("F33").Value = "Shreyas Iyer (c)"
$ws.Range("G33").Value = "88"
$ws.Range("H33").Value = "38"
$ws.Range("I33").Value = "7"
$ws.Range("J33").Value = "6"
$ws.Range("K33").Value = "231.57"

# Row 34
$ws.Range("A34:K34").NumberFormat = "@"
$ws.Range("A34").Value = " Dubai (DSC)"
$ws.Range("B34").Value = " September 25 2020"
$ws.Range("C34").Value = "Capitals won by 44 runs"
$ws.Range("D34").Value = "Delhi Capitals"
$ws.Range("E34").Value = "Chennai Super Kings"
$ws.Range("F34").Value = "Shreyas Iyer (c)"
$ws.Range("G34").Value = "26"
$ws.Range("H34").Value = "22"
$ws.Range("I34").Value = "1"
$ws.Range("J34").Value = "0"
$ws.Range("K34").Value = "118.18"

# Row 35
$ws.Range("A35:K35").NumberFormat = "@"
$ws.Range("A35").Value = " Sharjah"
$ws.Range("B35").Value = " October 09 2020"
$ws.Range("C35").Value = "Capitals won by 46 runs"
$ws.Range("D35").Value = "Delhi Capitals"
$ws.Range("E35").Value = "Rajasthan Royals"
$ws.Range("F35").Value = "Shreyas Iyer (c)"
$ws.Range("G35").Value = "22"
$ws.Range("H35").Value = "18"
$ws.Range("I35").Value = "4"
$ws.Range("J35").Value = "0"
$ws.Range("K35").Value = "122.22"
